$wb = $excel.ActiveWorkbook

# Insert a new worksheet "6.5" right before the existing "6.6" sheet,
# mirroring the workbook.xml <sheets> change (new sheetId=6, rId5; 6.6 -> rId6).
$sheet66 = $wb.Worksheets.Item("6.6")
$tp = $wb.Worksheets.Add($sheet66)
$tp.Name = "6.5"

# Insulation-level / thermal class data table for the "Transformador de
# Potencial" (Voltage Transformer) continuous-duty table.
$data = @(
  @(0.3,0.4,0.5,0.6,0.7,0.8,1,"Regime Contínuo"),
  @(60,50,50,50,40,40,30,20),
  @(110,90,70,70,70,60,60,40),
  @(180,150,120,120,110,100,80,60),
  @(310,260,200,200,180,160,140,100),
  @(530,450,340,340,300,270,250,150),
  @(890,750,570,570,500,500,430,230),
  @(1470,1240,1000,1000,900,850,740,370),
  @(2480,2060,1700,1700,1500,1400,1400,580),
  @(3300,2800,2000,2000,1900,1800,1500,930),
  @(5600,4700,3600,3600,3400,3000,1700,1500),
  @(9000,7600,5900,5900,5300,5000,4500,2400),
  @(13300,11600,9400,9400,8600,8000,7900,3700),
  @(17500,15700,13900,13900,13000,13000,13800,5900),
  @(26000,24000,21300,21300,21000,20000,24000,9300)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $tp.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Widen column H so the "Regime Contínuo" header is readable.
$tp.Columns.Item(8).ColumnWidth = 17.3

# Leave the selection on the last populated cell, matching the saved view.
$tp.Range("A15").Select() | Out-Null
